$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Faculty form creation completed:
# Replace the student email in D2 with the faculty's email address.
# (The cell keeps its existing hyperlink formatting; only the displayed
# text / underlying shared-string value changes.)
$ws.Range("D2").Value = "mannmehta2199@gmail.com"

# Move the active selection to D2 (matches the saved cursor position).
$ws.Range("D2").Select()
